# "colocando seção para baixar musica"
# Insert a new worksheet "musicas_de_sogro" between "lives" and "videos",
# give it a "link" header in A1, make it the active sheet, and (best
# effort) rename the built-in "Hyperlink" cell style to "Hiperlink".

$wb = $excel.ActiveWorkbook

$livesSheet = $wb.Worksheets.Item("lives")

# New sheet goes right after "lives" (i.e. before "videos").
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $livesSheet)
$newSheet.Name = "musicas_de_sogro"

# Seed it with the same "link" header used on the other tabs.
$newSheet.Range("A1").Value = "link"

# Leave the selection on A2, right under the header, and make this the
# sheet that's active/shown when the workbook is opened.
$newSheet.Range("A2").Select() | Out-Null
$newSheet.Activate() | Out-Null

# Best-effort: rename the built-in "Hyperlink" cell style to "Hiperlink"
# (pt-BR spelling). Harmless if the host doesn't persist style renames.
try {
    $wb.Styles.Item("Hyperlink").Name = "Hiperlink"
} catch {
}
